# Add a new "Godrej Capital" interview-experience entry to the tracker.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (24 - Neosoft) gets its Result filled in as "cleared".
$ws.Range("D24").Value = "cleared"

# New row 25: date, company, interview questions (no result yet).
$ws.Range("A25").Value = 45959
$ws.Range("B25").Value = "Godrej Capital"
$ws.Range("C25").Value = "aggragation, composition, how to avoid deadlock, volatile , automic, hashmap, design principals, design patterns find employee with salary greater than 20k`ndependency injection vs ioc, which di best`ndirectives in angular, how to optimize angular application"

# Leave the selection where the user would naturally end up after typing the new row.
$ws.Range("A26").Select() | Out-Null

Write-Output "Added Godrej Capital row"
